$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 (shifts old rows 7-17 down to 8-18)
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,R are constant across this block of rows,
# copy the same values used by the surrounding rows.
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = 'Vega Modelo de Temuco'
$ws.Range("C7").Value = 'La Araucanía'
$ws.Range("D7").Value = 44414
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 'Fruta'
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = 'Otros'
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = 'Chirimoya'
$ws.Range("K7").Value = 'Cultivar IV Región'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 3500
$ws.Range("O7").Value = 3500
$ws.Range("P7").Value = 3500
$ws.Range("Q7").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R7").Value = 'Provincia del Elquí'
$ws.Range("S7").Value = 3500
$ws.Range("T7").Value = 1

# Ensure the date cell keeps the date number format used by the rest of
# column D (style index 2 in the original workbook).
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
